# Apply historical-data metadata updates to the "attribute" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attribute")

# date_released: observed minimum date moves earlier (historical data now included)
$ws.Range("L4").Value = "2003-12-03"

# time_released: observed minimum/maximum times now reflect historical data
$ws.Range("L5").Value = "00:07:00"
$ws.Range("M5").Value = "24:00:00"

# release_site domain description: add the "P4" release location and reorder levels
$ws.Range("B3").Value = 'Location of release site. Levels = c("VB", "CCRB", "P4", NA)'

# number_released: minimum observed count drops
$ws.Range("L6").Value = 8

# median_fork_length_released: minimum observed length drops slightly
$ws.Range("L7").Value = 32

# day_or_night_release domain description: add "unknown" level
$ws.Range("B10").Value = 'Whether the release was conducted at day or night, levels = c("night", "day", NA, "unknown")'

# release_temp: min/max now reflect historical data
$ws.Range("L11").Value = 37.7
$ws.Range("M11").Value = 62.3

# release_flow: min/max now reflect historical data
$ws.Range("L12").Value = 194
$ws.Range("M12").Value = 1650

# release_turbidity: min/max now reflect historical data
$ws.Range("L13").Value = 0.58
$ws.Range("M13").Value = 93.8

# origin_released domain description: add "natural" level
$ws.Range("B14").Value = 'Origin of released fish. Levels = c(NA, "natural", "hatchery")'

# Widen the attribute_definition column to fit the longer domain text and
# update the active selection, matching the author's final workbook state.
$ws.Columns.Item(2).ColumnWidth = 68.66
$ws.Range("B16").Select() | Out-Null
